$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates (row 1) ---
# Rename existing headers
$ws.Range("U1").Value() = "GA Image 1"
$ws.Range("V1").Value() = "ACO Image 1"

# Create new header cells by copying formatting (style) from U1, then set their text
$ws.Range("U1").Copy($ws.Range("W1"))
$ws.Range("U1").Copy($ws.Range("Y1"))
$ws.Range("V1").Copy($ws.Range("X1"))
$ws.Range("V1").Copy($ws.Range("Z1"))
$ws.Range("W1").Value() = "GA Image 2"
$ws.Range("X1").Value() = "ACO Image 2"
$ws.Range("Y1").Value() = "GA Image 3"
$ws.Range("Z1").Value() = "ACO Image 3"

# --- Row 2 ---
$ws.Range("J2").Value() = 102.1070184418289
$ws.Range("K2").Value() = 104.0130965649983
$ws.Range("M2").Value() = 104.0130965649983
$ws.Range("N2").Value() = 104.0130965649983
$ws.Range("O2").Value() = 0.001
$ws.Range("P2").Value() = 0.0008
$ws.Range("Q2").Value() = 0.0009
$ws.Range("S2").Value() = 0.0008
$ws.Range("T2").Value() = 0.0008
$ws.Range("U2").Value() = "./imageResult/t5_1_GA_10.png"
$ws.Range("V2").Value() = "./imageResult/t5_1_ACO_10.png"
$ws.Range("W2").Value() = "./imageResult/t5_2_GA_10.png"
$ws.Range("X2").Value() = "./imageResult/t5_2_ACO_10.png"
$ws.Range("Y2").Value() = "./imageResult/t5_2_GA_10.png"
$ws.Range("Z2").Value() = "./imageResult/t5_2_ACO_10.png"

# --- Row 3 ---
$ws.Range("J3").Value() = 102.1070184418289
$ws.Range("K3").Value() = 102.1070184418289
$ws.Range("M3").Value() = 104.0130965649983
$ws.Range("N3").Value() = 104.0130965649983
$ws.Range("O3").Value() = 0.0009
$ws.Range("P3").Value() = 0.0009
$ws.Range("Q3").Value() = 0.0009
$ws.Range("S3").Value() = 0.0009
$ws.Range("T3").Value() = 0.0009
$ws.Range("U3").Value() = "./imageResult/t5_1_GA_50.png"
$ws.Range("V3").Value() = "./imageResult/t5_1_ACO_50.png"
$ws.Range("W3").Value() = "./imageResult/t5_2_GA_50.png"
$ws.Range("X3").Value() = "./imageResult/t5_2_ACO_50.png"
$ws.Range("Y3").Value() = "./imageResult/t5_2_GA_50.png"
$ws.Range("Z3").Value() = "./imageResult/t5_2_ACO_50.png"

# --- Row 4 ---
$ws.Range("J4").Value() = 102.1070184418289
$ws.Range("K4").Value() = 102.1070184418289
$ws.Range("M4").Value() = 104.0130965649983
$ws.Range("N4").Value() = 104.0130965649983
$ws.Range("O4").Value() = 0.0011
$ws.Range("P4").Value() = 0.0011
$ws.Range("Q4").Value() = 0.0011
$ws.Range("S4").Value() = 0.001
$ws.Range("T4").Value() = 0.001
$ws.Range("U4").Value() = "./imageResult/t5_1_GA_100.png"
$ws.Range("V4").Value() = "./imageResult/t5_1_ACO_100.png"
$ws.Range("W4").Value() = "./imageResult/t5_2_GA_100.png"
$ws.Range("X4").Value() = "./imageResult/t5_2_ACO_100.png"
$ws.Range("Y4").Value() = "./imageResult/t5_2_GA_100.png"
$ws.Range("Z4").Value() = "./imageResult/t5_2_ACO_100.png"

# --- Row 5 ---
$ws.Range("I5").Value() = 49.42117558688741
$ws.Range("J5").Value() = 46.26525516106483
$ws.Range("K5").Value() = 46.81591122239789
$ws.Range("M5").Value() = 31.22691510942754
$ws.Range("N5").Value() = 31.88252949105588
$ws.Range("O5").Value() = 0.0011
$ws.Range("P5").Value() = 0.0011
$ws.Range("Q5").Value() = 0.0011
$ws.Range("S5").Value() = 0.0009
$ws.Range("T5").Value() = 0.001
$ws.Range("U5").Value() = "./imageResult/burma14_1_GA_10.png"
$ws.Range("V5").Value() = "./imageResult/burma14_1_ACO_10.png"
$ws.Range("W5").Value() = "./imageResult/burma14_2_GA_10.png"
$ws.Range("X5").Value() = "./imageResult/burma14_2_ACO_10.png"
$ws.Range("Y5").Value() = "./imageResult/burma14_2_GA_10.png"
$ws.Range("Z5").Value() = "./imageResult/burma14_2_ACO_10.png"

# --- Row 6 ---
$ws.Range("I6").Value() = 37.67518177401836
$ws.Range("J6").Value() = 37.07000077260059
$ws.Range("K6").Value() = 37.08847053051498
$ws.Range("L6").Value() = 31.22691510942754
$ws.Range("M6").Value() = 31.22691510942754
$ws.Range("N6").Value() = 31.88252949105588
$ws.Range("O6").Value() = 0.0022
$ws.Range("P6").Value() = 0.0023
$ws.Range("Q6").Value() = 0.0022
$ws.Range("S6").Value() = 0.0016
$ws.Range("T6").Value() = 0.0018
$ws.Range("U6").Value() = "./imageResult/burma14_1_GA_50.png"
$ws.Range("V6").Value() = "./imageResult/burma14_1_ACO_50.png"
$ws.Range("W6").Value() = "./imageResult/burma14_2_GA_50.png"
$ws.Range("X6").Value() = "./imageResult/burma14_2_ACO_50.png"
$ws.Range("Y6").Value() = "./imageResult/burma14_2_GA_50.png"
$ws.Range("Z6").Value() = "./imageResult/burma14_2_ACO_50.png"

# --- Row 7 ---
$ws.Range("I7").Value() = 36.0214184683452
$ws.Range("J7").Value() = 34.86174060408727
$ws.Range("K7").Value() = 33.46723944553786
$ws.Range("M7").Value() = 31.45623383762054
$ws.Range("N7").Value() = 31.88252949105588
$ws.Range("O7").Value() = 0.0035
$ws.Range("P7").Value() = 0.0038
$ws.Range("Q7").Value() = 0.0034
$ws.Range("R7").Value() = 0.0024
$ws.Range("S7").Value() = 0.0028
$ws.Range("T7").Value() = 0.0025
$ws.Range("U7").Value() = "./imageResult/burma14_1_GA_100.png"
$ws.Range("V7").Value() = "./imageResult/burma14_1_ACO_100.png"
$ws.Range("W7").Value() = "./imageResult/burma14_2_GA_100.png"
$ws.Range("X7").Value() = "./imageResult/burma14_2_ACO_100.png"
$ws.Range("Y7").Value() = "./imageResult/burma14_2_GA_100.png"
$ws.Range("Z7").Value() = "./imageResult/burma14_2_ACO_100.png"

# --- Row 8 ---
$ws.Range("I8").Value() = 554427.2997867422
$ws.Range("J8").Value() = 568444.5940621259
$ws.Range("K8").Value() = 583017.4756071992
$ws.Range("L8").Value() = 49294.74163904427
$ws.Range("M8").Value() = 49215.61251916289
$ws.Range("N8").Value() = 49143.7729793856
$ws.Range("O8").Value() = 0.0087
$ws.Range("P8").Value() = 0.0085
$ws.Range("Q8").Value() = 0.0085
$ws.Range("R8").Value() = 0.0674
$ws.Range("S8").Value() = 0.0665
$ws.Range("T8").Value() = 0.0679
$ws.Range("U8").Value() = "./imageResult/lin318_1_GA_10.png"
$ws.Range("V8").Value() = "./imageResult/lin318_1_ACO_10.png"
$ws.Range("W8").Value() = "./imageResult/lin318_2_GA_10.png"
$ws.Range("X8").Value() = "./imageResult/lin318_2_ACO_10.png"
$ws.Range("Y8").Value() = "./imageResult/lin318_2_GA_10.png"
$ws.Range("Z8").Value() = "./imageResult/lin318_2_ACO_10.png"

# --- Row 9 ---
$ws.Range("I9").Value() = 511714.6596703269
$ws.Range("J9").Value() = 540311.5082501203
$ws.Range("K9").Value() = 533365.6925958826
$ws.Range("L9").Value() = 48563.49026440229
$ws.Range("M9").Value() = 48572.80330933771
$ws.Range("N9").Value() = 49215.61251916289
$ws.Range("O9").Value() = 0.0351
$ws.Range("P9").Value() = 0.0351
$ws.Range("Q9").Value() = 0.0351
$ws.Range("R9").Value() = 0.3182
$ws.Range("S9").Value() = 0.3186
$ws.Range("T9").Value() = 0.3203
$ws.Range("U9").Value() = "./imageResult/lin318_1_GA_50.png"
$ws.Range("V9").Value() = "./imageResult/lin318_1_ACO_50.png"
$ws.Range("W9").Value() = "./imageResult/lin318_2_GA_50.png"
$ws.Range("X9").Value() = "./imageResult/lin318_2_ACO_50.png"
$ws.Range("Y9").Value() = "./imageResult/lin318_2_GA_50.png"
$ws.Range("Z9").Value() = "./imageResult/lin318_2_ACO_50.png"

# --- Row 10 ---
$ws.Range("I10").Value() = 525354.1034295225
$ws.Range("J10").Value() = 503198.0617196271
$ws.Range("K10").Value() = 521017.7767133673
$ws.Range("L10").Value() = 48320.84193889733
$ws.Range("M10").Value() = 48835.97530486222
$ws.Range("N10").Value() = 48585.25741399533
$ws.Range("O10").Value() = 0.0701
$ws.Range("P10").Value() = 0.0687
$ws.Range("Q10").Value() = 0.0697
$ws.Range("R10").Value() = 0.6376
$ws.Range("S10").Value() = 0.6425
$ws.Range("T10").Value() = 0.6393
$ws.Range("U10").Value() = "./imageResult/lin318_1_GA_100.png"
$ws.Range("V10").Value() = "./imageResult/lin318_1_ACO_100.png"
$ws.Range("W10").Value() = "./imageResult/lin318_2_GA_100.png"
$ws.Range("X10").Value() = "./imageResult/lin318_2_ACO_100.png"
$ws.Range("Y10").Value() = "./imageResult/lin318_2_GA_100.png"
$ws.Range("Z10").Value() = "./imageResult/lin318_2_ACO_100.png"

